# VAN-1811: Prepare and write FUNCTIONAL test cases and test scripts
#
# The source data in Sheet1 has a single Ship-To customer (col K, shared
# string "Abiba Systems Private Limited") and a single "Previous Doc"
# number (col AX, 214002901789) repeated down every data row (rows 2-14).
# As part of preparing the functional-test fixture, these two values are
# replaced with anonymised placeholder values used throughout the sheet:
#   - col AX (Previous Doc)          -> "1663910091" (kept as literal text)
#   - col K  (Ship To Customer Name) -> "nwhjek726863"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column AX - Previous Doc. The replacement value is digits-only, so it is
# entered with a leading apostrophe to keep it stored as text (matching the
# rest of the anonymised data) instead of being reinterpreted as a number.
$ws.Range("AX2:AX14").Value = "'1663910091"

# Column K - Ship To Customer Name
$ws.Range("K2:K14").Value = "nwhjek726863"
